# Mark attendance for Session 18 (column X) as Absent ("A") for the
# students listed below. This mirrors marking those students absent
# for the 18th class session on the FA-II (Section A) attendance sheet.
#
# For each target row we copy the cell formatting from an existing
# "Absent" cell in the same row (so the style matches the rest of the
# sheet's "A" cells, border/alignment/etc. included) and then set the
# value to "A".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row number -> column letter of an existing "A" cell in that row whose
# formatting (style) should be copied onto column X.
$rows = [ordered]@{
    11 = "W"
    14 = "T"
    22 = "V"
    24 = "W"
    26 = "W"
    28 = "U"
    32 = "S"
    38 = "V"
    51 = "R"
    57 = "V"
    59 = "U"
    66 = "V"
    68 = "W"
    69 = "R"
    76 = "V"
    78 = "W"
}

foreach ($row in $rows.Keys) {
    $srcCol = $rows[$row]
    $srcCell = $ws.Range("$srcCol$row")
    $dstCell = $ws.Range("X$row")

    # Copy formatting from the existing "A" cell so the new cell matches
    # the styling already used for absences in the sheet.
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats

    # Mark the student absent for this session.
    $dstCell.Value = "A"
}

$excel.CutCopyMode = 0
